# Apply edits described by the commit:
#  - Mejorar generacion_excel: limpiar autor, formatear mes,
#    usar config.ini para validacion de codigos
#
# In practice this means, for the data rows of Sheet1:
#   * Column B (MES) is reformatted from "MM/YYYY" to the Spanish month name.
#   * Column A / I (Autor / nombre_periodista) drop the trailing ", Ver Biografía".
#   * Columns J, K, L, M, N (validated against config.ini) get re-extracted /
#     re-coded values.
#   * Columns O (TIEMPO_PROCESAMIENTO) and P (TIMESTAMP) reflect the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Juan Arias
$ws.Range("B2").Value = 'Enero'
$ws.Range("J2").Value = 'Sí, hombre'
$ws.Range("L2").Value = "'1"
$ws.Range("N2").Value = "'14"
$ws.Range("O2").Value = 8.25
$ws.Range("P2").Value = '2026-01-15 11:17:21'

# Row 3 - Arola Poch
$ws.Range("B3").Value = 'Febrero'
$ws.Range("J3").Value = 'Sí, hombre'
$ws.Range("K3").Value = "'1"
$ws.Range("L3").Value = "'4"
$ws.Range("M3").Value = "'4"
$ws.Range("N3").Value = "'1"
$ws.Range("O3").Value = 7.08
$ws.Range("P3").Value = '2026-01-15 11:17:29'

# Row 4 - María Porcel
$ws.Range("A4").Value = 'María Porcel'
$ws.Range("B4").Value = 'Enero'
$ws.Range("I4").Value = 'María Porcel'
$ws.Range("J4").Value = "'2"
$ws.Range("L4").Value = "'1"
$ws.Range("N4").Value = "'1"
$ws.Range("O4").Value = 5.87
$ws.Range("P4").Value = '2026-01-15 11:17:35'

# Row 5 - Daniel Soufi
$ws.Range("A5").Value = 'Daniel Soufi'
$ws.Range("B5").Value = 'Enero'
$ws.Range("I5").Value = 'Daniel Soufi'
$ws.Range("J5").Value = 'Sí, hombre'
$ws.Range("K5").Value = "'1"
$ws.Range("L5").Value = "'1"
$ws.Range("M5").Value = "'2"
$ws.Range("N5").Value = "'10"
$ws.Range("O5").Value = 6.48
$ws.Range("P5").Value = '2026-01-15 11:17:42'

# Row 6 - Pablo G. Bejerano
$ws.Range("B6").Value = 'Enero'
$ws.Range("J6").Value = 'Sí, hombre'
$ws.Range("L6").Value = "'1"
$ws.Range("N6").Value = "'15"
$ws.Range("O6").Value = 6.65
$ws.Range("P6").Value = '2026-01-15 11:17:49'

# Row 7 - Isabel Rubio
$ws.Range("A7").Value = 'Isabel Rubio'
$ws.Range("B7").Value = 'Enero'
$ws.Range("I7").Value = 'Isabel Rubio'
$ws.Range("J7").Value = 'Sí'
$ws.Range("N7").Value = "'15"
$ws.Range("O7").Value = 5.17
$ws.Range("P7").Value = '2026-01-15 11:17:55'

# Row 8 - David Trueba
$ws.Range("B8").Value = 'Abril'
$ws.Range("J8").Value = "'2"
$ws.Range("L8").Value = "'1"
$ws.Range("M8").Value = "'2"
$ws.Range("O8").Value = 5.8
$ws.Range("P8").Value = '2026-01-15 11:18:02'

# Row 9 - María De Los Ángeles Serrano Moral
$ws.Range("B9").Value = 'Junio'
$ws.Range("J9").Value = "'2"
$ws.Range("L9").Value = "'2"
$ws.Range("M9").Value = "'2"
$ws.Range("O9").Value = 4.95
$ws.Range("P9").Value = '2026-01-15 11:18:09'

# Row 10 - Borja Adsuara Varela
$ws.Range("B10").Value = 'Enero'
$ws.Range("J10").Value = "'2"
$ws.Range("L10").Value = "'2"
$ws.Range("N10").Value = "'15"
$ws.Range("O10").Value = 5.82
$ws.Range("P10").Value = '2026-01-15 11:18:16'

# Row 11 - Juan José Millás
$ws.Range("A11").Value = 'Juan José Millás'
$ws.Range("B11").Value = 'Junio'
$ws.Range("I11").Value = 'Juan José Millás'
$ws.Range("J11").Value = "'2"
$ws.Range("K11").Value = "'1"
$ws.Range("L11").Value = "'1"
$ws.Range("M11").Value = "'4"
$ws.Range("N11").Value = "'15"
$ws.Range("O11").Value = 5.26
$ws.Range("P11").Value = '2026-01-15 11:18:25'
